# Apply cryptos list update (prices / volume changes) per commit diff.
# Values are forced to Text (leading apostrophe) so numeric-looking
# strings like "1.00" or "3.189.17" are preserved exactly as text,
# then the cell Style is reset to "Normal" to drop the quote-prefix flag
# Excel applies when a value is entered with a leading apostrophe.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''87.634.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.06%  '
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = '''3.175.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -3.72%  '
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = '''  +0.00%  '
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = '''207.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -3.48%  '
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = '''610.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -3.21%  '
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = '''0.386'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -2.39%  '
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = '''0.672'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +4.12%  '
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = '''1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.01%  '
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = '''3.169.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -3.76%  '
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = '''0.539'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -8.69%  '
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = '''  -0.19%  '
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = '''0.0000244'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -8.65%  '
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = '''WrappedliquidstakedEther2.0'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = '''https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = '''3.770.47'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -3.41%  '
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = '''Toncoin'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = '''5.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -0.67%  '
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = '''WrappedBTC'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = '''https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = '''87.607.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.53%  '
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = '''Avalanche'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = '''https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = '''32.33'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -6.86%  '
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = '''3.183.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -3.51%  '
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = '''3.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +7.79%  '
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = '''13.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -5.76%  '
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = '''413.17'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -6.21%  '
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = '''8.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -8.49%  '
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = '''5.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -6.56%  '
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = '''5.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -0.03%  '
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = '''12.24'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -1.08%  '
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = '''3.349.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -3.46%  '
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = '''0.0000132'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.03%  '
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = '''73.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -5.07%  '
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = '''1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +0.06%  '
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = '''0.162'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -11.83%  '
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = '''1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.50%  '
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = '''544.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -2.28%  '
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = '''8.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -8.98%  '
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = '''1.32'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -9.86%  '
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = '''6.90'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -2.26%  '
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = '''1.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -6.62%  '
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = '''0.131'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -6.81%  '
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = '''21.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -4.45%  '
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = '''  +0.29%  '
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = '''  -0.12%  '
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = '''3.05'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +1.06%  '
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = '''  -0.01%  '
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = '''1.92'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -6.47%  '
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = '''0.374'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -7.80%  '
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = '''148.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -3.45%  '
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = '''173.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -4.34%  '
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = '''43.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -4.53%  '
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = '''0.125'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +3.20%  '
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = '''1.23'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -9.41%  '
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = '''3.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -7.53%  '
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = '''23.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -3.10%  '
$ws.Range("E51").Style = "Normal"
